# Updated symbol list on Sat Dec 31 16:59:29 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = "'246.82"
$ws.Range('D3').Value = "'26.49"
$ws.Range('D5').Value = "'0.05616"
$ws.Range('D6').Value = "'6.493"
$ws.Range('D7').Value = "'0.8138"
$ws.Range('D8').Value = "'0.8443"
$ws.Range('B9').Value = 'One'
$ws.Range('C9').Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D9').Value = "'0.009864"
$ws.Range('E9').Value = '8OneONEBestin24h'
$ws.Range('B10').Value = 'BitrueCoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D10').Value = "'0.02873"
$ws.Range('E10').Value = '9BitrueCoinBTR'
$ws.Range('B11').Value = 'BitMartToken'
$ws.Range('C11').Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D11').Value = "'0.09411"
$ws.Range('E11').Value = '10BitMartTokenBMX'
$ws.Range('B12').Value = 'BitForexToken'
$ws.Range('C12').Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D12').Value = "'0.001528"
$ws.Range('E12').Value = '11BitForexTokenBF'
$ws.Range('B13').Value = 'TigerCash'
$ws.Range('C13').Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D13').Value = "'0.006202"
$ws.Range('E13').Value = '12TigerCashTCH'
$ws.Range('B14').Value = 'LEO'
$ws.Range('C14').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D14').Value = "'3.593"
$ws.Range('E14').Value = '13LEOLEO'
$ws.Range('B15').Value = 'GateToken'
$ws.Range('C15').Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range('D15').Value = "'3.010"
$ws.Range('E15').Value = '14GateTokenGT'
$ws.Range('B16').Value = 'BTSEToken'
$ws.Range('C16').Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D16').Value = "'2.118"
$ws.Range('E16').Value = '15BTSETokenBTSE'
$ws.Range('B17').Value = 'BitpandaEcosystemToken'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D17').Value = "'0.3157"
$ws.Range('E17').Value = '16BitpandaEcosystemTokenBEST'
$ws.Range('B18').Value = 'WazirX'
$ws.Range('C18').Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range('D18').Value = "'0.1346"
$ws.Range('E18').Value = '17WazirXWRX'
$ws.Range('B19').Value = 'MandalaExchangeToken'
$ws.Range('C19').Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range('D19').Value = "'0.06988"
$ws.Range('E19').Value = '18MandalaExchangeTokenMDX'
$ws.Range('B20').Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').Value = "'0.03185"
$ws.Range('E20').Value = '19LiechtensteinCryptoassetsExchangeLCX'
$ws.Range('D22').Value = "'3.752"
$ws.Range('D23').Value = "'0.04667"
$ws.Range('D24').Value = "'0.1349"
$ws.Range('D25').Value = "'0.001249"
$ws.Range('D26').Value = "'0.004602"
$ws.Range('D27').Value = "'0.00009594"
$ws.Range('D28').Value = "'0.0001937"
$ws.Range('B41').Value = 'KickToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').Value = "'0.006156"
$ws.Range('E41').Value = '40KickTokenKICK'
$ws.Range('B42').Value = 'BKEXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').Value = "'0.1061"
$ws.Range('E42').Value = '41BKEXTokenBKK'
$ws.Range('D43').Value = "'0.002499"
$ws.Range('D44').Value = "'0.008906"
$ws.Range('D45').Value = "'0.00005290"
$ws.Range('D47').Value = "'0.1498"
$ws.Range('D48').Value = "'0.002524"
$ws.Range('D49').Value = "'0.00002099"
$ws.Range('D50').Value = "'0.0001999"
